$wb = $excel.ActiveWorkbook

# --- Repayment schedule: add column O (mirrors column N) for rows 2-8 ---
$repay = $wb.Worksheets.Item("Repayment schedule")
for ($r = 2; $r -le 8; $r++) {
    $repay.Range("N$r").Copy($repay.Range("O$r"))
}

# --- Summary: update selection ---
$summary = $wb.Worksheets.Item("Summary")
[void]$summary.Activate()
[void]$summary.Range("A4").Select()

# --- Edit Repayment Schedule: update selection (stays the active/selected tab) ---
$editSheet = $wb.Worksheets.Item("Edit Repayment Schedule")
[void]$editSheet.Activate()
[void]$editSheet.Range("B7").Select()
